$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14 (shifts rows 14-30 down to 15-31),
# matching the commit: "agregue 10 quincena al monto de 10000"
# (added the 10-fortnight term for the 10000 loan amount).
$ws.Rows("14:14").Insert()

# Fill in the new row's values: Monto=10000, Quincenas=10, Cuota=1500
$ws.Range("A14").Value = 10000
$ws.Range("B14").Value = 10
$ws.Range("C14").Value = 1500

# Match the author's final selection/view state in the saved file.
$ws.Range("B12").Select() | Out-Null
